# Insert a new row after row 176, pushing rows 177-293 down to 178-294,
# then populate the newly inserted row 177 with new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 177 (this shifts existing row 177 and below down by one)
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with the new data record
$ws.Cells.Item(177, 1).Value = 4
$ws.Cells.Item(177, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(177, 3).Value = "Los Lagos"
$ws.Cells.Item(177, 4).Value = Get-Date -Year 2022 -Month 7 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(177, 5).Value = 10
$ws.Cells.Item(177, 6).Value = "Fruta"
$ws.Cells.Item(177, 7).Value = 100104
$ws.Cells.Item(177, 8).Value = "Frutos de pepita"
$ws.Cells.Item(177, 9).Value = 100104005
$ws.Cells.Item(177, 10).Value = "Pera"
$ws.Cells.Item(177, 11).Value = "Packham's Triumph"
$ws.Cells.Item(177, 12).Value = "Primera"
$ws.Cells.Item(177, 13).Value = 150
$ws.Cells.Item(177, 14).Value = 16000
$ws.Cells.Item(177, 15).Value = 16000
$ws.Cells.Item(177, 16).Value = 16000
$ws.Cells.Item(177, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(177, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(177, 19).Value = 1067
$ws.Cells.Item(177, 20).Value = 15
